# Week 17 data logging + tiebreak-related numeric updates.
# Appends this week's per-game numbers to the running space-separated
# logs on the YDS and ST sheets, and bumps the season-total numeric
# cells on OFF / DEF / ST / TURNS / PEN accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet: append this week's R (rush) / P (pass) yardage-per-play
# logs for OFF (col B) and DEF (col C).
# ---------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Text + " 55 2 0 1 -1 14 2 1 4 12 1 9 0 -1 0 1 0 0 11 7 2 4 22 -1 5"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Text + " 6 3 4 6 6 -1 7 5 4 -5 5 -2 5 -4 1 2 3 11 4 0 2"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Text + " 9 9 13 5 9 12 5 9 24 11 24 19 2 23 11 21 14 4"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Text + " 1 4 9 15 7 4 31 21 3 2 21 14 6 24 14 9 8 20 14 1 8 32 9 21 3 12 8 3 6 6 27 10 33"

# ---------------------------------------------------------------
# ST sheet: append this week's per-kick logs (col B / col D rows 4-6).
# ---------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B4").Value = $stWs.Range("B4").Text + " 62"
$stWs.Range("B5").Value = $stWs.Range("B5").Text + " 7"
$stWs.Range("B6").Value = $stWs.Range("B6").Text + " 26 27 13 25 28"
$stWs.Range("D3").Value = $stWs.Range("D3").Text + " 38 32 46 38"
$stWs.Range("D4").Value = $stWs.Range("D4").Text + " 0 0 3 0"
$stWs.Range("D5").Value = $stWs.Range("D5").Text + " 0 0 0"

# ---------------------------------------------------------------
# OFF sheet: season totals through Week 17.
# ---------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("B2").Value = 6
$offWs.Range("C2").Value = 373
$offWs.Range("D2").Value = 19
$offWs.Range("E2").Value = 12
$offWs.Range("F2").Value = 101
$offWs.Range("G2").Value = 125
$offWs.Range("J2").Value = 46
$offWs.Range("L2").Value = 575
$offWs.Range("M2").Value = 340
$offWs.Range("O2").Value = 45
$offWs.Range("Q2").Value = 1044

$offWs.Range("C3").Value = 349
$offWs.Range("D3").Value = 3
$offWs.Range("E3").Value = 63
$offWs.Range("F3").Value = 207
$offWs.Range("I3").Value = 141
$offWs.Range("J3").Value = 115
$offWs.Range("N3").Value = 43

# ---------------------------------------------------------------
# DEF sheet: season totals through Week 17.
# ---------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 6
$defWs.Range("C2").Value = 450
$defWs.Range("F2").Value = 129
$defWs.Range("G2").Value = 120
$defWs.Range("J2").Value = 62
$defWs.Range("L2").Value = 618
$defWs.Range("M2").Value = 424
$defWs.Range("O2").Value = 32
$defWs.Range("P2").Value = 19
$defWs.Range("Q2").Value = 1143

$defWs.Range("C3").Value = 348
$defWs.Range("E3").Value = 75
$defWs.Range("F3").Value = 217
$defWs.Range("G3").Value = 77
$defWs.Range("H3").Value = 50
$defWs.Range("I3").Value = 122
$defWs.Range("J3").Value = 107
$defWs.Range("N3").Value = 39

# ---------------------------------------------------------------
# ST sheet: season totals through Week 17.
# ---------------------------------------------------------------
$stWs.Range("B2").Value = 127
$stWs.Range("D2").Value = 140
$stWs.Range("F2").Value = 25
$stWs.Range("G2").Value = 21
$stWs.Range("N2").Value = 4
$stWs.Range("O2").Value = 1

$stWs.Range("B3").Value = 78

# ---------------------------------------------------------------
# TURNS sheet: season totals through Week 17.
# ---------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 10
$turnsWs.Range("E3").Value = 24

# ---------------------------------------------------------------
# PEN sheet: season totals through Week 17.
# ---------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 25
